$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "42.889.43"
Set-TextValue "E2" "  -1.58%  "
Set-TextValue "D3" "2.560.33"
Set-TextValue "E3" "  -0.85%  "
Set-TextValue "E4" "  +0.12%  "
Set-TextValue "D5" "302.05"
Set-TextValue "E5" "  +0.55%  "
Set-TextValue "D6" "92.46"
Set-TextValue "E6" "  -3.56%  "
Set-TextValue "E7" "  -0.26%  "
Set-TextValue "E8" "  -0.05%  "
Set-TextValue "D9" "0.545"
Set-TextValue "E9" "  -1.51%  "
Set-TextValue "D10" "36.12"
Set-TextValue "E10" "  -2.51%  "
Set-TextValue "D11" "0.0810"
Set-TextValue "E11" "  -0.09%  "
Set-TextValue "D12" "7.77"
Set-TextValue "E12" "  -0.26%  "
Set-TextValue "E13" "  +6.84%  "
Set-TextValue "D14" "2.530.42"
Set-TextValue "E14" "  -2.21%  "
Set-TextValue "E15" "  -0.30%  "
Set-TextValue "D16" "14.20"
Set-TextValue "E16" "  -0.61%  "
Set-TextValue "D17" "42.939.74"
Set-TextValue "E17" "  -1.52%  "
Set-TextValue "D18" "0.0₃0993"
Set-TextValue "E18" "  +1.90%  "
Set-TextValue "D19" "12.63"
Set-TextValue "E19" "  +2.85%  "
Set-TextValue "E20" "  -0.56%  "
Set-TextValue "D21" "71.54"
Set-TextValue "E21" "  -1.85%  "
Set-TextValue "D22" "253.07"
Set-TextValue "E22" "  -4.37%  "
Set-TextValue "D23" "2.93"
Set-TextValue "E23" "  +0.43%  "
Set-TextValue "D24" "2.12"
Set-TextValue "E24" "  -4.43%  "
Set-TextValue "D25" "28.67"
Set-TextValue "E25" "  -2.11%  "
Set-TextValue "E26" "  -0.25%  "
Set-TextValue "D27" "10.25"
Set-TextValue "E27" "  +0.35%  "
Set-TextValue "D28" "37.13"
Set-TextValue "E28" "  -1.60%  "
Set-TextValue "E29" "  -3.89%  "
Set-TextValue "D30" "6.02"
Set-TextValue "E30" "  +1.03%  "
Set-TextValue "D31" "153.72"
Set-TextValue "E31" "  +1.42%  "
Set-TextValue "E32" "  -1.43%  "
Set-TextValue "E33" "  -5.86%  "
Set-TextValue "E34" "  -3.54%  "
Set-TextValue "D36" "18.11"
Set-TextValue "E36" "  +7.62%  "
Set-TextValue "D37" "0.113"
Set-TextValue "E37" "  -3.13%  "
Set-TextValue "E38" "  -0.41%  "
Set-TextValue "D39" "23.17"
Set-TextValue "E39" "  -3.76%  "
Set-TextValue "B40" "ApeXProtocol"
Set-TextValue "C40" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D40" "2.09"
Set-TextValue "E40" "  +30.38%  "
Set-TextValue "B41" "NEARProtocol"
Set-TextValue "C41" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D41" "3.42"
Set-TextValue "E41" "  -1.36%  "
Set-TextValue "B42" "VeChain"
Set-TextValue "C42" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D42" "0.0310"
Set-TextValue "E42" "  -0.92%  "
Set-TextValue "B43" "RenderToken"
Set-TextValue "C43" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D43" "3.87"
Set-TextValue "E43" "  +0.63%  "
Set-TextValue "B44" "Maker"
Set-TextValue "C44" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D44" "2.094.67"
Set-TextValue "E44" "  +0.89%  "
Set-TextValue "E45" "  +0.05%  "
Set-TextValue "E46" "  +1.07%  "
Set-TextValue "D47" "85.18"
Set-TextValue "E47" "  -2.94%  "
Set-TextValue "E48" "  +10.41%  "
Set-TextValue "D49" "106.70"
Set-TextValue "E49" "  +1.33%  "
Set-TextValue "D50" "2.809.74"
Set-TextValue "E50" "  -0.75%  "
Set-TextValue "B51" "Algorand"
Set-TextValue "C51" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D51" "0.191"
Set-TextValue "E51" "  +0.53%  "
